# "Termine la parte de gestion de perfiles y elimine la opcion de crear patentes"
#
# Appends the last three rows of the ID/NOMBRE/TIPO lookup table (columns
# G:I) on Sheet1, finishing the Encargado / Cocinero profiles and the
# GestionarComanda permission:
#   19 Encargado        G
#   20 Cocinero         G
#   21 GestionarComanda P

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("G21").Value = 19
$ws.Range("H21").Value = "Encargado"
$ws.Range("I21").Value = "G"

$ws.Range("G22").Value = 20
$ws.Range("H22").Value = "Cocinero"
$ws.Range("I22").Value = "G"

$ws.Range("G23").Value = 21
$ws.Range("H23").Value = "GestionarComanda"
$ws.Range("I23").Value = "P"

# Give the new NOMBRE cells the same yellow highlight + thin left/right
# borders used for the other entries in that column, then fan that exact
# look out to the other two new rows.
$h21 = $ws.Range("H21")
$h21.Interior.Color = 65535
$h21.Borders.Item(7).LineStyle = 1
$h21.Borders.Item(7).Weight = 2
$h21.Borders.Item(10).LineStyle = 1
$h21.Borders.Item(10).Weight = 2

$h21.Copy()
$ws.Range("H22:H23").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Leave the view where the author ended up while finishing this table.
$ws.Activate() | Out-Null
$ws.Range("A9").Select() | Out-Null
$excel.ActiveWindow.ScrollRow = 9
$ws.Range("G24").Select() | Out-Null
